$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, copying the format from the existing
# header cell G1 so it reuses the same style (bold, bordered, centered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the numeric values for the new "Save" column in the data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
